# Add a "cent_diac" column (inserted before the old "sign"/"comment" columns)
# and append three new word rows (climate crisis, refugee crisis, blockchain),
# renaming the old "degree centrality" header to "cent_total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D - this shifts the existing D ("sign") and E ("comment")
# columns one slot to the right (-> E and F) and keeps all other data intact.
$ws.Columns("D").Insert()

# New data rows appended below the existing table.
$ws.Range("A11").Value = "climate crisis"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "word of the year"

$ws.Range("A12").Value = "refugee crisis"

# Relabel the C1 header ("degree centrality" -> "cent_total") and give the
# freshly inserted column D a header of its own ("cent_diac").
$ws.Range("C1").Value = "cent_total"
$ws.Range("D1").Value = "cent_diac"

$ws.Range("A13").Value = "blockchain"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 5

# Match the saved selection state.
$ws.Range("D14").Select()
